$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MonitoringTools")

# Row 5: replace the old truncated "Demo on ElasticSearch..." stub text with the
# new "Logstash pitfalls to avoid those..." entry, and grow the row to fit it.
$ws.Range("A5").Value = "Logstash pitfalls to avoid those:`nhttp://logz.io/blog/5-logstash-pitfalls-and-how-to-avoid-them/`n"
$ws.Rows.Item(5).RowHeight = 72

# Row 6 (new): "How preprocess logs with logstash" entry.
$ws.Range("A6").Value = "How preprocess logs with logstash:`nhttp://blog.mmlac.com/how-to-pre-process-logs-with-logstash/"
$ws.Rows.Item(6).RowHeight = 43.2

# Match the author's final selection left on the sheet (cell A5).
$ws.Activate()
$ws.Range("A5").Select()
